$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Shift the old "footnote" row (old row 79) down two rows, to rows 81,
# so the new data rows can take rows 78-80. Insert() duplicates the
# formatting of the row being pushed down onto the two new blank rows.
$ws.Rows("79:80").Insert()

# Fill in the newly inserted + repurposed rows with the new daily figures.
$ws.Range("A78").Value = 43933
$ws.Range("B78").Value = 697
$ws.Range("C78").Value = 20958
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 4872

$ws.Range("A79").Value = 43934
$ws.Range("B79").Value = 1047
$ws.Range("C79").Value = 22005
$ws.Range("D79").Value = 195
$ws.Range("E79").Value = 5067

$ws.Range("A80").Value = 43935
$ws.Range("B80").Value = 770
$ws.Range("C80").Value = 22775
$ws.Range("D80").Value = 165
$ws.Range("E80").Value = 5232

# Extend the print area by two rows to keep pace with the new data.
$n = $wb.Names.Item("相談件数!Print_Area")
$n.RefersTo = "=相談件数!`$A`$1:`$E`$85"

# Match the author's final selection / frozen-pane placement.
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("C72").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("E80").Select()
